{"js": "// The document has a stray \"_GoBack\" bookmark sitting alone in its own\n// paragraph partway through the doc, and a trailing \"Author: ...\" byline\n// (preceded by a \"Spacing\"-styled blank paragraph) after the final \"JJ\"\n// signature. This change removes the byline block and relocates the\n// \"_GoBack\" bookmark to the end of the final \"JJ\" paragraph (its natural,\n// Word-maintained \"last edit\" position).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraph that consists solely of the \"_GoBack\" bookmark\n// (empty text) and the final \"JJ\" signature paragraph (the last paragraph\n// anywhere in the body whose text is exactly \"JJ\").\nlet lastJJIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"JJ\") {\n    lastJJIndex = i;\n  }\n}\n\nif (lastJJIndex === -1) {\n  throw new Error('Could not find a \"JJ\" paragraph to anchor the bookmark to.');\n}\n\n// Remove the old \"_GoBack\" bookmark from wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-add it at the end of the final \"JJ\" paragraph.\nconst jjParagraph = items[lastJJIndex];\nconst jjEnd = jjParagraph.getRange(\"End\");\njjEnd.insertBookmark(\"_GoBack\");\n\n// Delete the trailing \"Spacing\" blank paragraph and the \"Author: ...\"\n// byline paragraph that follow the final \"JJ\" paragraph.\nconst spacingParagraph = items[lastJJIndex + 1];\nconst authorParagraph = items[lastJJIndex + 2];\nspacingParagraph.load(\"text\");\nauthorParagraph.load(\"text\");\nawait context.sync();\n\nif (!/^Author:/.test(authorParagraph.text)) {\n  throw new Error(\"Unexpected paragraph after the trailing Spacing paragraph; aborting to avoid deleting the wrong content.\");\n}\n\nspacingParagraph.delete();\nawait context.sync();\nauthorParagraph.delete();\n\nawait context.sync();\n", "ps1": "# The document has a stray \"_GoBack\" bookmark sitting alone in its own\n# paragraph partway through the doc, and a trailing \"Author: ...\" byline\n# (preceded by a \"Spacing\"-styled blank paragraph) after the final \"JJ\"\n# signature. This change removes the byline block and relocates the\n# \"_GoBack\" bookmark to the end of the final \"JJ\" paragraph (its natural,\n# Word-maintained \"last edit\" position).\n\n$d = $word.ActiveDocument\n\n# Step 1: drop the old \"_GoBack\" bookmark from wherever it currently sits.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Step 2: locate the final paragraph whose only content is \"JJ\" (the\n# trailing signature at the very end of the document).\n$jjParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"JJ\") {\n        $jjParaIndex = $i\n    }\n}\nif ($jjParaIndex -eq -1) {\n    throw \"Could not find a 'JJ' paragraph to anchor the bookmark to.\"\n}\n$jjPara = $d.Paragraphs.Item($jjParaIndex)\n\n# Step 3: re-insert \"_GoBack\" at the end of that paragraph's text (right\n# after \"JJ\", before the paragraph mark). A bare zero-length range dropped\n# exactly on that boundary does not reliably anchor a bookmark, so a\n# one-character placeholder is inserted, wrapped in the bookmark, and then\n# cleared back out -- leaving the bookmark correctly collapsed in place.\n$insertionPoint = $jjPara.Range.Duplicate\n$insertionPoint.Collapse(0)\n$insertionPoint.MoveEnd(1, -1)\n$insertionPoint.InsertAfter([char]7)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n$insertionPoint.Text = \"\"\n\n# Step 4: delete the trailing \"Spacing\" blank paragraph and the\n# \"Author: ...\" byline paragraph that follow the final \"JJ\" paragraph.\n# Paragraph objects/indices are re-fetched after each deletion rather than\n# reused, since a stale reference into a shifted collection is a no-op.\n$spacingPara = $d.Paragraphs.Item($jjParaIndex + 1)\nif ($spacingPara.Style.NameLocal -ne \"Spacing\") {\n    throw \"Unexpected paragraph after the final 'JJ' paragraph; aborting to avoid deleting the wrong content.\"\n}\n$authorTextPreview = $d.Paragraphs.Item($jjParaIndex + 2).Range.Text\nif ($authorTextPreview -notmatch \"^Author:\") {\n    throw \"Unexpected paragraph after the trailing Spacing paragraph; aborting to avoid deleting the wrong content.\"\n}\n\n$spacingPara.Range.Delete()\n\n$authorPara = $d.Paragraphs.Item($jjParaIndex + 1)\n$authorPara.Range.Delete()\n"}
